# "adding support for port ranges"
#
# MasterProtected (sheet 1): the Port value on row 2 becomes a range
# (22,  3389 -> 22-100,  3389).
#
# FirewallRulesToValidate (sheet 2):
#   - "Azure Non prod" -> "Azure Non prod VM range" and
#     "IOD- Non-Prod"  -> "IOD- Non-Prod VM" for the existing Azure-VM-range
#     rows (2,3,5,6), plus a previously-missing row 4 that gets the same
#     correct Source/Destination description (row 4's Destination
#     Description had mistakenly held "Azure" instead of "IOD- Non-Prod VM").
#   - Row 2's Port becomes a range ("22,  443" -> "22-25,  443").
#   - Four new rows (7-10) describing the 10.108.0.1 http rule, with ports
#     80, 101, the range "80-105" and the range "101-105" (single ports and
#     ranges side-by-side to validate the new range support).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # MasterProtected
$ws2 = $wb.Worksheets.Item(2)   # FirewallRulesToValidate

# ---------------------------------------------------------------------------
# MasterProtected
# ---------------------------------------------------------------------------
$ws1.Range("G2").Value = "22-100,  3389"
$ws1.Columns.Item(7).ColumnWidth = 14
$ws1.Range("C29").Select() | Out-Null

# ---------------------------------------------------------------------------
# FirewallRulesToValidate
# ---------------------------------------------------------------------------

# Row 2 - Azure Non prod VM range / IOD- Non-Prod VM, port becomes a range
$ws2.Range("A2").Value = "Azure Non prod VM range"
$ws2.Range("C2").Value = "IOD- Non-Prod VM"
$ws2.Range("F2").Value = "22-25,  443"

# Row 3 - same source/destination description rename
$ws2.Range("A3").Value = "Azure Non prod VM range"
$ws2.Range("C3").Value = "IOD- Non-Prod VM"

# Row 4 - rename + fix the destination description (was wrongly "Azure")
$ws2.Range("A4").Value = "Azure Non prod VM range"
$ws2.Range("B4").Value = "10.200.0.1"
$ws2.Range("C4").Value = "IOD- Non-Prod VM"
$ws2.Range("D4").Value = "10.201.0.1"
$ws2.Range("E4").Value = "TCP"
$ws2.Range("F4").Value = 3389
$ws2.Range("G4").Value = "SSH"

# Row 5 - same source/destination description rename
$ws2.Range("A5").Value = "Azure Non prod VM range"
$ws2.Range("C5").Value = "IOD- Non-Prod VM"

# Row 6 - same source/destination description rename
$ws2.Range("A6").Value = "Azure Non prod VM range"
$ws2.Range("C6").Value = "IOD- Non-Prod VM"

# Rows 7-10 - new rule entries for the 10.108.0.1 http range, showing both
# single ports (80, 101) and port ranges (80-105, 101-105)
$ws2.Range("A7").Value = "Azure Non prod"
$ws2.Range("B7").Value = "10.108.0.1"
$ws2.Range("B7").WrapText = $true
$ws2.Range("C7").Value = "IOD- Non-Prod"
$ws2.Range("D7").Value = "20.0.0.0/16"
$ws2.Range("E7").Value = "udp"
$ws2.Range("F7").Value = 80
$ws2.Range("G7").Value = "http"

$ws2.Range("A8").Value = "Azure Non prod"
$ws2.Range("B8").Value = "10.108.0.1"
$ws2.Range("B8").WrapText = $true
$ws2.Range("C8").Value = "IOD- Non-Prod"
$ws2.Range("D8").Value = "20.0.0.0/16"
$ws2.Range("E8").Value = "udp"
$ws2.Range("F8").Value = 101
$ws2.Range("G8").Value = "http"

$ws2.Range("A9").Value = "Azure Non prod"
$ws2.Range("B9").Value = "10.108.0.1"
$ws2.Range("B9").WrapText = $true
$ws2.Range("C9").Value = "IOD- Non-Prod"
$ws2.Range("D9").Value = "20.0.0.0/16"
$ws2.Range("E9").Value = "udp"
$ws2.Range("F9").Value = "80-105"
$ws2.Range("G9").Value = "http"

$ws2.Range("A10").Value = "Azure Non prod"
$ws2.Range("B10").Value = "10.108.0.1"
$ws2.Range("B10").WrapText = $true
$ws2.Range("C10").Value = "IOD- Non-Prod"
$ws2.Range("D10").Value = "20.0.0.0/16"
$ws2.Range("E10").Value = "udp"
$ws2.Range("F10").Value = "101-105"
$ws2.Range("G10").Value = "http"

# Column-width touch-ups that came along with the new, wider content
$ws2.Columns.Item(1).ColumnWidth = 24.67
$ws2.Columns.Item(3).ColumnWidth = 16.67
$ws2.Columns.Item(6).ColumnWidth = 18.17

$ws2.Activate() | Out-Null
$ws2.Range("F11").Select() | Out-Null
